# Append four new species-observation records (rows 57-60) to the
# "Artfynd" sheet, mirroring the existing rows already in the table.
#
# Columns I (Antal), Y/AA (Start/Slutdatum) hold digit- or date-shaped
# text ("14", "2023-09-17", ...) that must stay TEXT, not be coerced into
# numbers/dates by Excel's normal type-sniffing on assignment. We force
# that by switching the cell to a text NumberFormat before assigning the
# value, then clearing the format again so no stray formatting is left
# behind on the cell (matching the plain, unstyled source cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Column layout (same order as the sheet header):
# A Id | B Taxonsorteringsordning | C Valideringsstatus | D Rödlistade
# E TaxonId | F Artnamn | G Vetenskapligt namn | H Auktor | I Antal
# J Enhet | K Ålder-Stadium | L Kön | M Aktivitet | N Metod | P Lokalnamn
# Q Ost | R Nord | S Noggrannhet | T Län | U Kommun | V Provins
# W Församling | Y Startdatum | Z Starttid | AA Slutdatum | AB Sluttid
# AD Ej återfunnen | AE Osäker artbestämning | AG Ospontan
# AT Bestämningsår | AW Rapportör | AX Observatörer | AY Projektnamn
$rows = @(
    [ordered]@{
        Row = 57
        A = 112145545; B = 96348; C = "Ovaliderad"; D = "VU"; E = 220787
        F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
        I = "14"; J = "plantor/tuvor"; K = ""; L = ""; N = ""
        P = "Bennarby, Upl"; Q = 653038.3046146344; R = 6675340.776511455; S = 4
        T = "Uppsala"; U = "Östhammar"; V = "Uppland"; W = "Dannemora"
        Y = "2023-09-17"; Z = "11:07"; AA = "2023-09-17"; AB = "11:08"
        AD = $false; AE = $false; AG = $false; AT = ""
        AW = "Annika Rastén"; AX = "Annika Rastén"; AY = ""
    },
    [ordered]@{
        Row = 58
        A = 112145544; B = 96348; C = "Ovaliderad"; D = "VU"; E = 220787
        F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
        I = "17"; J = "plantor/tuvor"; K = ""; L = ""; N = ""
        P = "Bennarby, Upl"; Q = 653023.8828454959; R = 6675363.578861667; S = 4
        T = "Uppsala"; U = "Östhammar"; V = "Uppland"; W = "Dannemora"
        Y = "2023-09-17"; Z = "11:03"; AA = "2023-09-17"; AB = "11:04"
        AD = $false; AE = $false; AG = $false; AT = ""
        AW = "Annika Rastén"; AX = "Annika Rastén"; AY = ""
    },
    [ordered]@{
        Row = 59
        A = 112145539; B = 90018; C = "Ovaliderad"; D = "LC"; E = 1339
        F = "Brandticka"; G = "Pycnoporellus fulgens"; H = "(Fr.) Donk"
        I = ""; J = ""; K = ""; N = ""
        P = "Bennarby, Upl"; Q = 652996.6865235955; R = 6675310.152517678; S = 4
        T = "Uppsala"; U = "Östhammar"; V = "Uppland"; W = "Dannemora"
        Y = "2023-09-17"; Z = "10:44"; AA = "2023-09-17"; AB = "10:44"
        AD = $false; AE = $false; AG = $false; AT = ""
        AW = "Annika Rastén"; AX = "Annika Rastén"; AY = ""
    },
    [ordered]@{
        Row = 60
        A = 112145535; B = 56543; C = "Ovaliderad"; D = "NT"; E = 103021
        F = "Talltita"; G = "Poecile montanus"; H = "(Conrad von Baldenstein, 1827)"
        I = "1"; K = ""; L = ""; M = "permanent revir"; N = ""
        P = "Smigruvan, Upl"; Q = 653011.7621751076; R = 6675152.417146614; S = 84
        T = "Uppsala"; U = "Östhammar"; V = "Uppland"; W = "Dannemora"
        Y = "2023-09-17"; Z = "10:10"; AA = "2023-09-17"; AB = "10:10"
        AD = $false; AE = $false; AG = $false; AT = ""
        AW = "Annika Rastén"; AX = "Annika Rastén"; AY = ""
    }
)

# Columns whose source values are digit-/date-shaped text that must be
# pinned to Text so Excel doesn't reinterpret them as numbers or dates.
$textForceCols = @("I", "Y", "AA")

foreach ($rec in $rows) {
    $r = $rec.Row
    foreach ($col in $rec.Keys) {
        if ($col -eq "Row") { continue }
        $value = $rec[$col]
        $cell = $ws.Range($col + $r)
        if ($value -is [bool]) {
            $cell.Value = $value
        } elseif (($textForceCols -contains $col) -and ($value -ne "")) {
            Set-TextValue $cell $value
        } else {
            $cell.Value = $value
        }
    }
}

"Added rows 57-60 to " + $ws.Name
